$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.930.47"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").Value = "3.065.23"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "520.45"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").Value = "134.98"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.064.89"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "0.454"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "3.594.75"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "25.10"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "56.967.62"
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("D18").Value = "3.068.62"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "5.85"
$ws.Range("E19").Value = "  -4.63%  "
$ws.Range("D20").Value = "12.39"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "347.20"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "68.51"
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("E26").Value = "  -2.31%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "0.0₃0865"
$ws.Range("E28").Value = "  -6.91%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("D31").Value = "1.85"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").Value = "5.82"
$ws.Range("E32").Value = "  -8.80%  "
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").Value = "159.22"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("E36").Value = "  -5.59%  "
$ws.Range("D37").Value = "5.97"
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("D38").Value = "25.24"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").Value = "1.22"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Value = "0.0655"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("D42").Value = "4.02"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "0.689"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "2.387.45"
$ws.Range("E44").Value = "  +5.13%  "
$ws.Range("D45").Value = "36.60"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "3.102.63"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "0.0260"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("D49").Value = "0.948"
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("D50").Value = "5.94"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("D51").Value = "19.54"
$ws.Range("E51").Value = "  -5.31%  "
